$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 190
$ws.Range("F5").Value = 887
$ws.Range("F7").Value = 365
$ws.Range("F9").Value = 171
$ws.Range("F12").Value = 166
$ws.Range("C13").Value = "北京·AINI二次元派对【免票展会】"
$ws.Range("D13").Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Range("E13").Value = "2024.08.10 10:00-08.10 16:00"
$ws.Range("F13").Value = 792
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=89601"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202407/eIryW6Up1721208870214.jpeg"
$ws.Range("C14").Value = "北京·GOJO超次元动漫游戏嘉年华15th"
$ws.Range("D14").Value = "小关路39号 北投购物公园"
$ws.Range("E14").Value = "2024.08.10 09:20-08.11 17:00"
$ws.Range("F14").Value = 4192
$ws.Range("G14").Value = 6.6
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=85223"
$ws.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202407/DlsfmegR1720613188484.jpeg"
$ws.Range("F16").Value = 6034
$ws.Range("C17").Value = "北京·狐妖小红娘专题聚会【免票活动】"
$ws.Range("E17").Value = "2024.08.10 14:00-08.10 18:00"
$ws.Range("F17").Value = 58
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=90238"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202408/mL8ytYCG1722578125040.jpeg"
$ws.Range("F18").Value = 465
$ws.Range("F19").Value = 2323
$ws.Range("F21").Value = 465
$ws.Range("F22").Value = 9094
$ws.Range("F23").Value = 148
$ws.Range("F24").Value = 2447
$ws.Range("F25").Value = 2303
$ws.Range("F26").Value = 1386
$ws.Range("F27").Value = 235
$ws.Range("F28").Value = 1953
$ws.Range("F30").Value = 57
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 49
$ws.Range("F36").Value = 90
$ws.Range("F37").Value = 1214
$ws.Range("F38").Value = 72
$ws.Range("F40").Value = 234
$ws.Range("F41").Value = 1522
$ws.Range("F42").Value = 2487
$ws.Range("F43").Value = 917
$ws.Range("F47").Value = 28
$ws.Range("F49").Value = 15
$ws.Range("F50").Value = 62
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 190
$ws.Range("F5").Value = 365
$ws.Range("F6").Value = 171
$ws.Range("F7").Value = 166
$ws.Range("C8").Value = "北京·AINI二次元派对【免票展会】"
$ws.Range("D8").Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Range("E8").Value = "2024.08.10 10:00-08.10 16:00"
$ws.Range("F8").Value = 792
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89601"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202407/eIryW6Up1721208870214.jpeg"
$ws.Range("C9").Value = "北京·GOJO超次元动漫游戏嘉年华15th"
$ws.Range("D9").Value = "小关路39号 北投购物公园"
$ws.Range("E9").Value = "2024.08.10 09:20-08.11 17:00"
$ws.Range("F9").Value = 4191
$ws.Range("G9").Value = 6.6
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85223"
$ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202407/DlsfmegR1720613188484.jpeg"
$ws.Range("C10").Value = "北京·LookLook剧情式沉浸游戏互动动漫嘉年华（取消）"
$ws.Range("D10").Value = "东村文化创意产业园A1-2 五道杠实景片场"
$ws.Range("E10").Value = "2024.08.10 09:30-08.11 17:30"
$ws.Range("F10").Value = 1025
$ws.Range("G10").Value = "不可售"
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=84741"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202405/WH4KKudj1716880619473.jpeg"
$ws.Range("F13").Value = 6034
$ws.Range("C14").Value = "北京·狐妖小红娘专题聚会【免票活动】"
$ws.Range("E14").Value = "2024.08.10 14:00-08.10 18:00"
$ws.Range("F14").Value = 58
$ws.Range("G14").Value = 58
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=90238"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202408/mL8ytYCG1722578125040.jpeg"
$ws.Range("F15").Value = 465
$ws.Range("F16").Value = 2323
$ws.Range("F19").Value = 465
$ws.Range("F20").Value = 9093
$ws.Range("F22").Value = 2447
$ws.Range("F24").Value = 2303
$ws.Range("F25").Value = 2424
$ws.Range("F26").Value = 1386
$ws.Range("F27").Value = 235
$ws.Range("F28").Value = 1953
$ws.Range("F30").Value = 57
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 49
$ws.Range("F37").Value = 86
$ws.Range("F39").Value = 1214
$ws.Range("F40").Value = 72
$ws.Range("F42").Value = 234
$ws.Range("F43").Value = 1522
$ws.Range("F44").Value = 2487
$ws.Range("F45").Value = 917
$ws.Range("F48").Value = 15
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 887
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 148
$ws.Range("F21").Value = 28
$ws.Range("F22").Value = 62
$ws.Range("F23").Value = 62

Write-Output "Applied all changes"